$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Insert new slide 13 "Remaining Tasks" (pushes old slide13.. down by one)
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Add(13, 2)

$s1.Shapes.Item(1).TextFrame.TextRange.Text = "Remaining Tasks"

$body1 = $s1.Shapes.Item(2).TextFrame.TextRange
$lines1 = @("", "CURRENT - Authentication", "Data Management functional", "Updated UI", "Updated Viewer", "Google Cardboard ", "")

# First paragraph already exists (empty) - format it, then append the rest one
# at a time so the ruler-driven indent/margin override lands on every
# paragraph (the host stamps marL/indent onto whichever paragraph is
# "paragraph 1" at the moment the ruler is touched).
$idx = 1
foreach ($ln in $lines1) {
    if ($idx -gt 1) {
        $body1.InsertAfter("`r" + $ln)
    }
    $para = $body1.Paragraphs($idx, 1)
    $para.IndentLevel = 2
    $para.ParagraphFormat.Bullet.Type = 0
    $ruler = $s1.Shapes.Item(2).TextFrame2.Ruler
    $ruler.Levels.Item(2).LeftMargin = 15.84
    $ruler.Levels.Item(2).FirstMargin = 0
    $idx++
}

# ---------------------------------------------------------------------------
# 2) Insert new slide 21 "Problem & Solution - Authentication" (between the
#    relocated "Winter Week 6" slide and the relocated "Problems/Solutions"
#    slide).
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Add(21, 2)

$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Problem & Solution - Authentication"

$body2 = $s2.Shapes.Item(2).TextFrame.TextRange
$body2.Text = "One of the biggest problems"
$body2.InsertAfter("`rScratch code -> adapting existing code")
$body2.InsertAfter("`rNot understanding code")
$body2.InsertAfter("`rContact with ")
$body2.InsertAfter("code creator")
